# Update the "Overview" sheet of the yearly report:
#  - shift the twelve-month reporting period columns one year forward
#    (drop the 1396/12 column, add a new 1401/12 column)
#  - refresh every numeric data cell (rows 10,13,15,16,17,19,20,26,27)
#    with the newly shifted-in figures, per the new read_price algorithm

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# ---- Header row 8 (G&A expense table) ----
$ws.Range("E8").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1401/12"

# ---- Header row 24 (personnel count table) ----
$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# ---- Data rows: column values shift left one year, newest year appended in column I ----
$data = @{
    10 = @(154878, 852188, 868460, 1595081, 2332128)
    13 = @(1726,   2071,   1282,   2546,    3254)
    15 = @(325,    450,    696,    855,     1025)
    16 = @(1596,   1793,   3137,   9412,    23129)
    17 = @(46615,  54642,  60752,  73531,   91869)
    19 = @(67916,  57951,  96190,  97878,   167758)
    20 = @(273056, 969095, 1030517, 1779303, 2619163)
    26 = @(264,    256,    234,    234,     469)
    27 = @(90,     50,     24,     24,      44)
}

$cols = @("E", "F", "G", "H", "I")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
